# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) across all job sheets with latest market-board figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5907.391
$ws.Range("I51").Value = 5316.8335
$ws.Range("J51").Value = 6115.8237
$ws.Range("K51").Value = 5316.8335
$ws.Range("L51").Value = 6115.8237
$ws.Range("M51").Value = -4832.8335
$ws.Range("N51").Value = -7083.8237

$ws.Range("H62").Value = 10774
$ws.Range("I62").Value = 14865.375
$ws.Range("J62").Value = 5738.4614
$ws.Range("K62").Value = 14865.375
$ws.Range("L62").Value = 5738.4614
$ws.Range("M62").Value = -14241.375
$ws.Range("N62").Value = -6986.4614

$ws.Range("H65").Value = 10774
$ws.Range("I65").Value = 14865.375
$ws.Range("J65").Value = 5738.4614
$ws.Range("K65").Value = 74326.875
$ws.Range("L65").Value = 28692.307
$ws.Range("M65").Value = -71206.875
$ws.Range("N65").Value = -34932.307

$ws.Range("H98").Value = 2952
$ws.Range("I98").Value = 1559.4482
$ws.Range("J98").Value = 8000
$ws.Range("K98").Value = 1559.4482
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = -61.44820000000004
$ws.Range("N98").Value = -10996

$ws.Range("H122").Value = 2952
$ws.Range("I122").Value = 1559.4482
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 4678.3446
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -2228.3446
$ws.Range("N122").Value = -28900

$ws.Range("H132").Value = 5188.5454
$ws.Range("I132").Value = 1894.1714
$ws.Range("J132").Value = 18000
$ws.Range("K132").Value = 5682.5142
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -3152.5142
$ws.Range("N132").Value = -59060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2335.2856
$ws.Range("I2").Value = 3358.5
$ws.Range("J2").Value = 971
$ws.Range("K2").Value = 3358.5
$ws.Range("L2").Value = 971
$ws.Range("M2").Value = -3245.5
$ws.Range("N2").Value = -1197

$ws.Range("H32").Value = 1369514.9
$ws.Range("I32").Value = 1518569.2
$ws.Range("J32").Value = 3183.1667
$ws.Range("K32").Value = 1518569.2
$ws.Range("L32").Value = 3183.1667
$ws.Range("M32").Value = -1518282.2
$ws.Range("N32").Value = -3757.1667

$ws.Range("H39").Value = 2000
$ws.Range("I39").Value = 2000
$ws.Range("K39").Value = 2000
$ws.Range("M39").Value = -1480

$ws.Range("H116").Value = 2335.2856
$ws.Range("I116").Value = 3358.5
$ws.Range("J116").Value = 971
$ws.Range("K116").Value = 3358.5
$ws.Range("L116").Value = 971
$ws.Range("M116").Value = -1064.5
$ws.Range("N116").Value = -5559

$ws.Range("H132").Value = 27770.15
$ws.Range("I132").Value = 54964.105
$ws.Range("J132").Value = 3166.0952
$ws.Range("K132").Value = 164892.315
$ws.Range("L132").Value = 9498.285600000001
$ws.Range("M132").Value = -162362.315
$ws.Range("N132").Value = -14558.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2335.2856
$ws.Range("I3").Value = 3358.5
$ws.Range("J3").Value = 971
$ws.Range("K3").Value = 3358.5
$ws.Range("L3").Value = 971
$ws.Range("M3").Value = -3244.5
$ws.Range("N3").Value = -1199

$ws.Range("H101").Value = 26000
$ws.Range("J101").Value = 26000
$ws.Range("L101").Value = 26000
$ws.Range("N101").Value = -32490

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 11724.875
$ws.Range("I94").Value = 2049.75
$ws.Range("J94").Value = 21400
$ws.Range("K94").Value = 2049.75
$ws.Range("L94").Value = 21400
$ws.Range("M94").Value = -1598.75
$ws.Range("N94").Value = -22302

$ws.Range("H104").Value = 20259
$ws.Range("I104").Value = 20259
$ws.Range("K104").Value = 20259
$ws.Range("M104").Value = -17638

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""

$ws.Range("H107").Value = 276.92682
$ws.Range("I107").Value = 261.03845
$ws.Range("J107").Value = 304.46667
$ws.Range("K107").Value = 261.03845
$ws.Range("L107").Value = 304.46667
$ws.Range("M107").Value = 1658.96155
$ws.Range("N107").Value = -4144.46667

$ws.Range("H132").Value = 17244342
$ws.Range("I132").Value = 26318076
$ws.Range("K132").Value = 78954228
$ws.Range("M132").Value = -78951698

$ws.Range("H134").Value = 12501986
$ws.Range("I134").Value = 17858726
$ws.Range("J134").Value = 2924.75
$ws.Range("K134").Value = 53576178
$ws.Range("L134").Value = 8774.25
$ws.Range("M134").Value = -53573643
$ws.Range("N134").Value = -13844.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 352.625
$ws.Range("I12").Value = 182.25
$ws.Range("J12").Value = 523
$ws.Range("K12").Value = 546.75
$ws.Range("L12").Value = 1569
$ws.Range("M12").Value = -373.75
$ws.Range("N12").Value = -1915

$ws.Range("H14").Value = 2086.8
$ws.Range("I14").Value = 2086.8
$ws.Range("K14").Value = 6260.400000000001
$ws.Range("M14").Value = -6087.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""

$ws.Range("H105").Value = 40671
$ws.Range("J105").Value = 40671
$ws.Range("L105").Value = 40671
$ws.Range("N105").Value = -47659

$ws.Range("H123").Value = 10271.667
$ws.Range("J123").Value = 10271.667
$ws.Range("L123").Value = 10271.667
$ws.Range("N123").Value = -15171.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 845.8913
$ws.Range("I46").Value = 681.8333
$ws.Range("J46").Value = 1436.5
$ws.Range("K46").Value = 681.8333
$ws.Range("L46").Value = 1436.5
$ws.Range("M46").Value = -493.8333
$ws.Range("N46").Value = -1812.5

$ws.Range("H101").Value = 8749.5
$ws.Range("J101").Value = 8749.5
$ws.Range("L101").Value = 8749.5
$ws.Range("N101").Value = -15239.5

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""

$ws.Range("H103").Value = 20999.666
$ws.Range("J103").Value = 20999.666
$ws.Range("L103").Value = 20999.666
$ws.Range("N103").Value = -23343.666

$ws.Range("H104").Value = 16548
$ws.Range("J104").Value = 16548
$ws.Range("L104").Value = 16548
$ws.Range("N104").Value = -23536

$ws.Range("H122").Value = 2750
$ws.Range("I122").Value = 2750
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8250
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5800
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 30600
$ws.Range("J101").Value = 30600
$ws.Range("L101").Value = 30600
$ws.Range("N101").Value = -37090

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""

$ws.Range("H132").Value = 2296.2563
$ws.Range("I132").Value = 1613.05
$ws.Range("J132").Value = 3015.4211
$ws.Range("K132").Value = 4839.15
$ws.Range("L132").Value = 9046.263300000001
$ws.Range("M132").Value = -2309.15
$ws.Range("N132").Value = -14106.2633
